$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) values - prefixed with apostrophe to force text
# and avoid Excel auto-converting numeric-looking strings to Number type,
# matching the original inlineStr/text cell type.
$ws.Range("D2").Value = "'26.491.70"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "'1.733.63"
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'246.53"
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "'0.4888"
$ws.Range("E7").Value = "  +1.73%  "
$ws.Range("D8").Value = "'0.2663"
$ws.Range("E8").Value = "  -0.33%  "
$ws.Range("D9").Value = "'0.06265"
$ws.Range("E9").Value = "  +0.46%  "
$ws.Range("D10").Value = "'1.729.43"
$ws.Range("E10").Value = "  -0.49%  "
$ws.Range("D11").Value = "'0.07016"
$ws.Range("E11").Value = "  -1.73%  "
$ws.Range("D12").Value = "'15.66"
$ws.Range("E12").Value = "  -0.76%  "
$ws.Range("D13").Value = "'4.591"
$ws.Range("E13").Value = "  +1.10%  "
$ws.Range("E14").Value = "  -1.83%  "
$ws.Range("D15").Value = "'77.36"
$ws.Range("E15").Value = "  +0.50%  "
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("D17").Value = "'0.000007491"
$ws.Range("E17").Value = "  +8.52%  "
$ws.Range("D18").Value = "'26.479.18"
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("D20").Value = "'11.50"
$ws.Range("E20").Value = "  -2.23%  "
$ws.Range("D21").Value = "'1.952.12"
$ws.Range("E21").Value = "  -0.48%  "
$ws.Range("D22").Value = "'4.567"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "'8.702"
$ws.Range("E24").Value = "  -2.27%  "
$ws.Range("D25").Value = "'140.57"
$ws.Range("E25").Value = "  +3.57%  "
$ws.Range("D26").Value = "'15.43"
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("E27").Value = "  -0.28%  "
$ws.Range("D28").Value = "'1.764"
$ws.Range("E28").Value = "  -2.20%  "
$ws.Range("D29").Value = "'107.79"
$ws.Range("E29").Value = "  +1.01%  "
$ws.Range("D30").Value = "'4.018"
$ws.Range("E30").Value = "  +0.58%  "
$ws.Range("D31").Value = "'0.08012"
$ws.Range("E31").Value = "  +1.47%  "
$ws.Range("D32").Value = "'3.696"
$ws.Range("E32").Value = "  -0.61%  "
$ws.Range("D33").Value = "'0.04568"
$ws.Range("E33").Value = "  -0.39%  "
$ws.Range("D34").Value = "'1.000"
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("D35").Value = "'2.612"
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("E36").Value = "  +0.69%  "
$ws.Range("D37").Value = "'0.6336"
$ws.Range("E37").Value = "  -0.27%  "
$ws.Range("D38").Value = "'0.8946"
$ws.Range("E38").Value = "  -3.72%  "
$ws.Range("D39").Value = "'2.008"
$ws.Range("E39").Value = "  +1.24%  "
$ws.Range("D40").Value = "'2.396"
$ws.Range("E40").Value = "  -2.05%  "
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("E42").Value = "  -0.67%  "
$ws.Range("D43").Value = "'101.60"
$ws.Range("E43").Value = "  -8.90%  "
$ws.Range("D44").Value = "'5.402"
$ws.Range("E44").Value = "  -5.46%  "
$ws.Range("D45").Value = "'0.3877"
$ws.Range("E45").Value = "  -0.71%  "
$ws.Range("D46").Value = "'6.914"
$ws.Range("E46").Value = "  -0.17%  "
$ws.Range("D47").Value = "'0.1184"
$ws.Range("E47").Value = "  -1.04%  "
$ws.Range("D48").Value = "'0.05392"
$ws.Range("E48").Value = "  +1.10%  "
$ws.Range("E49").Value = "  -0.96%  "
$ws.Range("D50").Value = "'30.44"
$ws.Range("E50").Value = "  -1.28%  "
$ws.Range("D51").Value = "'1.257"
$ws.Range("E51").Value = "  +0.27%  "
